$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E (price / volume) to remain text so values
# like "1.008" or "21.29" are not auto-converted to numbers by Excel,
# then clear the temporary formatting so the cell style index is
# unchanged (back to the default/general style).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.581.70"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "1.840.03"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -2.59%  "
$ws.Range("D5").Value = "316.54"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("D7").Value = "0.4304"
$ws.Range("E7").Value = "  -1.92%  "
$ws.Range("D8").Value = "0.3723"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "0.07279"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").Value = "0.8693"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("D11").Value = "21.29"
$ws.Range("E11").Value = "  -1.25%  "
$ws.Range("D12").Value = "1.848.61"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "6.716"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("D14").Value = "5.380"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "0.07101"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "88.61"
$ws.Range("E16").Value = "  +4.27%  "
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  -2.77%  "
$ws.Range("D18").Value = "0.000008955"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").Value = "15.30"
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("D21").Value = "27.581.07"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").Value = "5.174"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").Value = "10.98"
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").Value = "2.074.02"
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").Value = "1.982"
$ws.Range("E25").Value = "  -4.35%  "
$ws.Range("D26").Value = "153.99"
$ws.Range("E26").Value = "  -3.06%  "
$ws.Range("D27").Value = "18.47"
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").Value = "2.160"
$ws.Range("E28").Value = "  +8.66%  "
$ws.Range("D29").Value = "5.309"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "0.08890"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "0.7719"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").Value = "4.507"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").Value = "2.900"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("E36").Value = "  -2.48%  "
$ws.Range("D37").Value = "1.125"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("D38").Value = "0.01966"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "0.05291"
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.147"
$ws.Range("E40").Value = "  +4.18%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.879"
$ws.Range("E41").Value = "  +1.18%  "
$ws.Range("D42").Value = "0.5102"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").Value = "0.1680"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").Value = "8.739"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "10.62"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").Value = "106.73"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("D47").Value = "0.4730"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").Value = "0.06440"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "1.006"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("E51").Value = "  -2.53%  "

$ws.Range("D2:E51").ClearFormats()
